# test_traits_dupTraitNames.xlsx - remove the stray "unit1"/"unit2" sample
# values that were left in the Units column (M2:M3) of the Template sheet,
# and update the active selection to reflect that range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

# Clear M2 and M3 (previously "unit1" / "unit2"). This also drops those two
# now-unreferenced strings from the shared string table.
$ws.Range("M2").Value = $null
$ws.Range("M3").Value = $null

# Reflect the new selection/active cell on the sheet.
$ws.Range("M2:M3").Select()
